$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Catalogue")
$ws.Range("F1").Value = "Nombre de pièces ou poids du conditionnement"
$ws.Range("J5").Select()
